# Auto-generated Excel COM-interop edit script
# Applies updated Leve market-price metrics (columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 38.77778
$ws.Range("I8").Value = 29.294117
$ws.Range("K8").Value = 87.882351
$ws.Range("M8").Value = 51.117649
$ws.Range("L17").Value = 7501638
$ws.Range("H17").Value = 2000496.8
$ws.Range("J17").Value = 2500546
$ws.Range("N17").Value = -7501974
$ws.Range("I69").Value = 7751
$ws.Range("N69").Value = -26553.332
$ws.Range("K69").Value = 23253
$ws.Range("J69").Value = 8268.444
$ws.Range("M69").Value = -22379
$ws.Range("L69").Value = 24805.332
$ws.Range("H69").Value = 8174.364
$ws.Range("N72").Value = -83151.996
$ws.Range("M72").Value = -65391
$ws.Range("L72").Value = 74415.996
$ws.Range("I72").Value = 7751
$ws.Range("J72").Value = 8268.444
$ws.Range("H72").Value = 8174.364
$ws.Range("K72").Value = 69759
$ws.Range("I92").Value = 395.5
$ws.Range("J92").Value = 20181
$ws.Range("N92").Value = -22677
$ws.Range("L92").Value = 20181
$ws.Range("K92").Value = 395.5
$ws.Range("H92").Value = 6214.7646
$ws.Range("M92").Value = 852.5
$ws.Range("L100").Value = 2900
$ws.Range("H100").Value = 4058.6316
$ws.Range("J100").Value = 2900
$ws.Range("N100").Value = -3982
$ws.Range("H112").Value = 47589.59
$ws.Range("K112").Value = 1007800.02
$ws.Range("I112").Value = 335933.34
$ws.Range("J112").Value = 2061.6316
$ws.Range("M112").Value = -1006692.02
$ws.Range("L112").Value = 6184.8948
$ws.Range("N112").Value = -8400.8948
$ws.Range("M135").Value = -6294.782999999999
$ws.Range("J135").Value = 155
$ws.Range("H135").Value = 885.7692
$ws.Range("N135").Value = -6465
$ws.Range("I135").Value = 981.087
$ws.Range("K135").Value = 8829.782999999999
$ws.Range("L135").Value = 1395
$ws.Range("K137").Value = 4592.1084
$ws.Range("N137").Value = -12237.5712
$ws.Range("M137").Value = -2042.1084
$ws.Range("H137").Value = 2065.25
$ws.Range("J137").Value = 2379.1904
$ws.Range("I137").Value = 1530.7028
$ws.Range("L137").Value = 7137.5712

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K4").Value = 6516.3335
$ws.Range("M4").Value = -6400.3335
$ws.Range("H4").Value = 11789.8
$ws.Range("I4").Value = 6516.3335
$ws.Range("L32").Value = 3964.3333
$ws.Range("N32").Value = -4538.3333
$ws.Range("H32").Value = 158332.8
$ws.Range("K32").Value = 165924.69
$ws.Range("M32").Value = -165637.69
$ws.Range("J32").Value = 3964.3333
$ws.Range("I32").Value = 165924.69
$ws.Range("K45").Value = 1976.9286
$ws.Range("H45").Value = 2223.5625
$ws.Range("I45").Value = 1976.9286
$ws.Range("M45").Value = -1599.9286
$ws.Range("M63").Value = -111117978
$ws.Range("H63").Value = 66676780
$ws.Range("I63").Value = 111118664
$ws.Range("K63").Value = 111118664
$ws.Range("K66").Value = 555593320
$ws.Range("I66").Value = 111118664
$ws.Range("H66").Value = 66676780
$ws.Range("M66").Value = -555589888

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 9814.846
$ws.Range("M20").Value = -9567.846
$ws.Range("H20").Value = 7777.1665
$ws.Range("K20").Value = 9814.846
$ws.Range("I94").Value = 2876.5925
$ws.Range("H94").Value = 3146.1875
$ws.Range("K94").Value = 2876.5925
$ws.Range("M94").Value = -2425.5925
$ws.Range("I105").Value = 1949.4445
$ws.Range("K105").Value = 1949.4445
$ws.Range("M105").Value = -202.4445000000001
$ws.Range("H105").Value = 2177.8696
$ws.Range("H107").Value = 20094226
$ws.Range("K107").Value = 107018.95
$ws.Range("M107").Value = -105098.95
$ws.Range("I107").Value = 107018.95
$ws.Range("J135").Value = 59999
$ws.Range("H135").Value = 49999.5
$ws.Range("N135").Value = -70139
$ws.Range("L135").Value = 59999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L16").Value = 3666.3333
$ws.Range("M16").Value = -3710.1667
$ws.Range("N16").Value = -4240.3333
$ws.Range("J16").Value = 3666.3333
$ws.Range("K16").Value = 3997.1667
$ws.Range("H16").Value = 3886.889
$ws.Range("I16").Value = 3997.1667
$ws.Range("H22").Value = 1617.75
$ws.Range("I22").Value = 361.375
$ws.Range("M22").Value = -11.375
$ws.Range("K22").Value = 361.375
$ws.Range("I50").Value = 28000
$ws.Range("J50").Value = 57500
$ws.Range("L50").Value = 57500
$ws.Range("N50").Value = -58750
$ws.Range("H50").Value = 47666.668
$ws.Range("K50").Value = 28000
$ws.Range("M50").Value = -27375
$ws.Range("H74").Value = 39764.5
$ws.Range("L74").Value = 47908.57
$ws.Range("N74").Value = -49656.57
$ws.Range("J74").Value = 47908.57
$ws.Range("L77").Value = 143725.71
$ws.Range("H77").Value = 39764.5
$ws.Range("N77").Value = -152461.71
$ws.Range("J77").Value = 47908.57
$ws.Range("L105").Value = 1299
$ws.Range("H105").Value = 1331
$ws.Range("J105").Value = 1299
$ws.Range("N105").Value = -4793
$ws.Range("M113").Value = -1827.1667
$ws.Range("I113").Value = 3997.1667
$ws.Range("H113").Value = 3886.889
$ws.Range("L113").Value = 3666.3333
$ws.Range("N113").Value = -8006.3333
$ws.Range("K113").Value = 3997.1667
$ws.Range("J113").Value = 3666.3333
$ws.Range("N132").Value = -18034.25
$ws.Range("H132").Value = 3851.7273
$ws.Range("I132").Value = 3581.4285
$ws.Range("J132").Value = 4324.75
$ws.Range("K132").Value = 10744.2855
$ws.Range("L132").Value = 12974.25
$ws.Range("M132").Value = -8214.2855

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L2").Value = 552.999984
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 92.166664
$ws.Range("N2").Value = -778.999984
$ws.Range("K2").Value = 204
$ws.Range("M2").Value = -91
$ws.Range("L17").Value = 1365
$ws.Range("H17").Value = 455
$ws.Range("J17").Value = 455
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("I17").Value = 0
$ws.Range("N17").Value = -1703
$ws.Range("M38").Value = -98.63637999999997
$ws.Range("N38").Value = -802
$ws.Range("L38").Value = 108
$ws.Range("K38").Value = 445.63638
$ws.Range("H38").Value = 124.42857
$ws.Range("J38").Value = 36
$ws.Range("I38").Value = 148.54546
$ws.Range("K40").Value = 104.5
$ws.Range("I40").Value = 26.125
$ws.Range("H40").Value = 218.0625
$ws.Range("M40").Value = -35.5
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1185
$ws.Range("I70").Value = 500
$ws.Range("H70").Value = 500
$ws.Range("M73").Value = -408
$ws.Range("H73").Value = 500
$ws.Range("K73").Value = 1500
$ws.Range("I73").Value = 500

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N132").Value = -24521
$ws.Range("H132").Value = 560284.4399999999
$ws.Range("I132").Value = 837183.2
$ws.Range("J132").Value = 6487
$ws.Range("K132").Value = 2511549.6
$ws.Range("L132").Value = 19461
$ws.Range("M132").Value = -2509019.6

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M122").Value = -10702
$ws.Range("K122").Value = 13152
$ws.Range("H122").Value = 4624
$ws.Range("I122").Value = 4384
$ws.Range("N132").Value = -33720.911
$ws.Range("H132").Value = 7740
$ws.Range("I132").Value = 3750
$ws.Range("J132").Value = 9553.637000000001
$ws.Range("K132").Value = 11250
$ws.Range("L132").Value = 28660.911
$ws.Range("M132").Value = -8720

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 7940507.5
$ws.Range("L62").Value = 16593.334
$ws.Range("M62").Value = -7939883.5
$ws.Range("N62").Value = -17841.334
$ws.Range("H62").Value = 1337245.8
$ws.Range("K62").Value = 7940507.5
$ws.Range("J62").Value = 16593.334
$ws.Range("N65").Value = -89206.67
$ws.Range("M65").Value = -39699417.5
$ws.Range("I65").Value = 7940507.5
$ws.Range("H65").Value = 1337245.8
$ws.Range("J65").Value = 16593.334
$ws.Range("K65").Value = 39702537.5
$ws.Range("L65").Value = 82966.67
$ws.Range("N81").Value = -36378970
$ws.Range("H81").Value = 11769900
$ws.Range("J81").Value = 18188424
$ws.Range("L81").Value = 36376848
$ws.Range("L84").Value = 181884240
$ws.Range("J84").Value = 18188424
$ws.Range("H84").Value = 11769900
$ws.Range("N84").Value = -181894848
$ws.Range("M113").Value = -4787
$ws.Range("I113").Value = 2319
$ws.Range("H113").Value = 7182
$ws.Range("L113").Value = 40998
$ws.Range("N113").Value = -45338
$ws.Range("K113").Value = 6957
$ws.Range("J113").Value = 13666
$ws.Range("N132").Value = -20252.8568
$ws.Range("H132").Value = 1182409
$ws.Range("I132").Value = 2006550.4
$ws.Range("J132").Value = 5064.2856
$ws.Range("K132").Value = 6019651.199999999
$ws.Range("L132").Value = 15192.8568
$ws.Range("M132").Value = -6017121.199999999
$ws.Range("J136").Value = 30000
$ws.Range("M136").ClearContents()
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("N136").Value = -95100
$ws.Range("L136").Value = 90000
$ws.Range("H136").Value = 30000

